$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.710887666666667
$ws.Range("H2").Value = 11.132663
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 22.43285649763678
$ws.Range("R2").Value = 201.8957084787311
$ws.Range("S2").Value = 0.8160840232643366
$ws.Range("T2").Value = 0.8160840232643367

# Row 3
$ws.Range("G3").Value = 3.710887666666667
$ws.Range("H3").Value = 11.132663
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 2.532470548865556
$ws.Range("R3").Value = 22.79223493979
$ws.Range("S3").Value = 0.09212864864242169
$ws.Range("T3").Value = 0.09212864864242169

# Row 4
$ws.Range("G4").Value = 3.710887666666667
$ws.Range("H4").Value = 11.132663
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 2.523088187881667
$ws.Range("R4").Value = 22.707793690935
$ws.Range("S4").Value = 0.09178732809324164
$ws.Range("T4").Value = 0.09178732809324165
